$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Debit")
$ws.Range("B2").Value = "Wed Feb 08 19:06:01 EST 2023"
$ws.Range("B3").Value = "Wed Feb 08 19:06:15 EST 2023"
$ws.Range("B4").Value = "Wed Feb 08 19:06:29 EST 2023"
$ws.Range("B5").Value = "Wed Feb 08 19:06:43 EST 2023"
$ws.Range("B6").Value = "Wed Feb 08 19:06:57 EST 2023"
$ws.Range("B7").Value = "Wed Feb 08 19:07:11 EST 2023"
$ws.Range("B8").Value = "Wed Feb 08 19:07:25 EST 2023"

$ws = $wb.Worksheets.Item("Debit-ZeroDollar")
$ws.Range("B2").Value = "Wed Feb 08 19:07:40 EST 2023"
$ws.Range("B3").Value = "Wed Feb 08 19:07:55 EST 2023"
$ws.Range("B4").Value = "Wed Feb 08 19:08:08 EST 2023"
$ws.Range("B5").Value = "Wed Feb 08 19:08:22 EST 2023"
$ws.Range("B6").Value = "Wed Feb 08 19:08:36 EST 2023"
$ws.Range("B7").Value = "Wed Feb 08 19:08:50 EST 2023"
$ws.Range("B8").Value = "Wed Feb 08 19:09:04 EST 2023"

$ws = $wb.Worksheets.Item("Debit-Void")
$ws.Range("B2").Value = "Wed Feb 08 19:09:18 EST 2023"
$ws.Range("B3").Value = "Wed Feb 08 19:09:42 EST 2023"
$ws.Range("B4").Value = "Wed Feb 08 19:10:05 EST 2023"
$ws.Range("B5").Value = "Wed Feb 08 19:10:29 EST 2023"
$ws.Range("B6").Value = "Wed Feb 08 19:10:52 EST 2023"
$ws.Range("B7").Value = "Wed Feb 08 19:11:16 EST 2023"
$ws.Range("B8").Value = "Wed Feb 08 19:11:39 EST 2023"

$ws = $wb.Worksheets.Item("Debit-Credit")
$ws.Range("B2").Value = "Wed Feb 08 19:12:04 EST 2023"
$ws.Range("B3").Value = "Wed Feb 08 19:12:29 EST 2023"
$ws.Range("B4").Value = "Wed Feb 08 19:12:54 EST 2023"
$ws.Range("B5").Value = "Wed Feb 08 19:13:18 EST 2023"
$ws.Range("B6").Value = "Wed Feb 08 19:13:43 EST 2023"
$ws.Range("B7").Value = "Wed Feb 08 19:14:08 EST 2023"
$ws.Range("B8").Value = "Wed Feb 08 19:14:32 EST 2023"

$ws = $wb.Worksheets.Item("Debit-Credit-Void")
$ws.Range("B2").Value = "Wed Feb 08 19:14:58 EST 2023"
$ws.Range("B3").Value = "Wed Feb 08 19:15:31 EST 2023"
$ws.Range("B4").Value = "Wed Feb 08 19:16:06 EST 2023"
$ws.Range("B5").Value = "Wed Feb 08 19:16:40 EST 2023"
$ws.Range("B6").Value = "Wed Feb 08 19:17:14 EST 2023"
$ws.Range("B7").Value = "Wed Feb 08 19:17:48 EST 2023"
$ws.Range("B8").Value = "Wed Feb 08 19:18:22 EST 2023"

$ws = $wb.Worksheets.Item("Debit-MRF")
$ws.Range("B2").Value = "Wed Feb 08 19:18:57 EST 2023"
$ws.Range("B3").Value = "Wed Feb 08 19:19:17 EST 2023"
$ws.Range("B4").Value = "Wed Feb 08 19:19:37 EST 2023"
$ws.Range("B5").Value = "Wed Feb 08 19:19:58 EST 2023"
$ws.Range("B6").Value = "Wed Feb 08 19:20:18 EST 2023"
$ws.Range("B7").Value = "Wed Feb 08 19:20:38 EST 2023"
$ws.Range("B8").Value = "Wed Feb 08 19:20:59 EST 2023"
$ws.Range("B9").Value = "Wed Feb 08 19:21:20 EST 2023"
$ws.Range("B10").Value = "Wed Feb 08 19:21:40 EST 2023"
$ws.Range("B11").Value = "Wed Feb 08 19:22:01 EST 2023"
$ws.Range("B12").Value = "Wed Feb 08 19:22:21 EST 2023"
$ws.Range("B13").Value = "Wed Feb 08 19:22:41 EST 2023"
$ws.Range("B14").Value = "Wed Feb 08 19:23:03 EST 2023"
$ws.Range("B15").Value = "Wed Feb 08 19:23:23 EST 2023"
$ws.Range("B16").Value = "Wed Feb 08 19:23:44 EST 2023"
$ws.Range("B17").Value = "Wed Feb 08 19:24:04 EST 2023"
$ws.Range("B18").Value = "Wed Feb 08 19:24:24 EST 2023"
$ws.Range("B19").Value = "Wed Feb 08 19:24:45 EST 2023"
$ws.Range("B20").Value = "Wed Feb 08 19:25:05 EST 2023"

$ws = $wb.Worksheets.Item("Void-MRF")
$ws.Range("B2").Value = "Wed Feb 08 19:25:26 EST 2023"
$ws.Range("B3").Value = "Wed Feb 08 19:25:38 EST 2023"
$ws.Range("B4").Value = "Wed Feb 08 19:25:49 EST 2023"
$ws.Range("B5").Value = "Wed Feb 08 19:26:00 EST 2023"
$ws.Range("B6").Value = "Wed Feb 08 19:26:12 EST 2023"
$ws.Range("A7").Value = "Fail"
$ws.Range("B7").Value = "Wed Feb 08 19:26:23 EST 2023"
$ws.Range("B8").Value = "Wed Feb 08 19:26:34 EST 2023"
$ws.Range("B9").Value = "Wed Feb 08 19:26:46 EST 2023"
$ws.Range("B10").Value = "Wed Feb 08 19:26:57 EST 2023"
$ws.Range("B11").Value = "Wed Feb 08 19:27:09 EST 2023"

$ws = $wb.Worksheets.Item("Credit-MRF")
$ws.Range("B2").Value = "Wed Feb 08 19:27:21 EST 2023"
$ws.Range("B3").Value = "Wed Feb 08 19:27:34 EST 2023"
$ws.Range("B4").Value = "Wed Feb 08 19:27:46 EST 2023"
$ws.Range("B5").Value = "Wed Feb 08 19:27:59 EST 2023"
$ws.Range("B6").Value = "Wed Feb 08 19:28:12 EST 2023"
$ws.Range("B7").Value = "Wed Feb 08 19:28:25 EST 2023"
$ws.Range("B8").Value = "Wed Feb 08 19:28:38 EST 2023"
$ws.Range("B9").Value = "Wed Feb 08 19:28:51 EST 2023"
$ws.Range("B10").Value = "Wed Feb 08 19:29:04 EST 2023"
$ws.Range("B11").Value = "Wed Feb 08 19:29:17 EST 2023"
$ws.Range("B12").Value = "Wed Feb 08 19:29:29 EST 2023"
$ws.Range("B13").Value = "Wed Feb 08 19:29:42 EST 2023"

$ws = $wb.Worksheets.Item("Debit-RemID-Pipe")
$ws.Range("B2").Value = "Wed Feb 08 19:32:48 EST 2023"
$ws.Range("B3").Value = "Wed Feb 08 19:33:02 EST 2023"
$ws.Range("B4").Value = "Wed Feb 08 19:33:17 EST 2023"
$ws.Range("B5").Value = "Wed Feb 08 19:33:31 EST 2023"
$ws.Range("B6").Value = "Wed Feb 08 19:33:45 EST 2023"
$ws.Range("B7").Value = "Wed Feb 08 19:33:59 EST 2023"
$ws.Range("B8").Value = "Wed Feb 08 19:34:13 EST 2023"

$ws = $wb.Worksheets.Item("DebitVoid-RemID-Pipe")
$ws.Range("B2").Value = "Wed Feb 08 19:34:28 EST 2023"
$ws.Range("B3").Value = "Wed Feb 08 19:34:52 EST 2023"
$ws.Range("B4").Value = "Wed Feb 08 19:35:15 EST 2023"
$ws.Range("B5").Value = "Wed Feb 08 19:35:39 EST 2023"
$ws.Range("B6").Value = "Wed Feb 08 19:36:03 EST 2023"
$ws.Range("B7").Value = "Wed Feb 08 19:36:26 EST 2023"
$ws.Range("B8").Value = "Wed Feb 08 19:36:50 EST 2023"

$ws = $wb.Worksheets.Item("DebitCredit-RemID-Pipe")
$ws.Range("B2").Value = "Wed Feb 08 19:29:56 EST 2023"
$ws.Range("B3").Value = "Wed Feb 08 19:30:20 EST 2023"
$ws.Range("B4").Value = "Wed Feb 08 19:30:45 EST 2023"
$ws.Range("B5").Value = "Wed Feb 08 19:31:09 EST 2023"
$ws.Range("B6").Value = "Wed Feb 08 19:31:34 EST 2023"
$ws.Range("B7").Value = "Wed Feb 08 19:31:59 EST 2023"
$ws.Range("B8").Value = "Wed Feb 08 19:32:23 EST 2023"
